$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue "D2" "27.581.05"
Set-TextValue "E2" "  +0.89%  "
Set-TextValue "D3" "1.641.57"
Set-TextValue "E3" "  -0.59%  "
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "212.42"
Set-TextValue "E5" "  -0.41%  "
Set-TextValue "D6" "0.537"
Set-TextValue "E6" "  +4.70%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.15%  "
Set-TextValue "D8" "22.88"
Set-TextValue "E8" "  -3.89%  "
Set-TextValue "E9" "  -1.52%  "
Set-TextValue "E10" "  -0.53%  "
Set-TextValue "E11" "  +1.44%  "
Set-TextValue "D12" "1.872.79"
Set-TextValue "E12" "  -0.72%  "
Set-TextValue "D13" "1.635.43"
Set-TextValue "E13" "  -1.56%  "
Set-TextValue "E14" "  -0.88%  "
Set-TextValue "D15" "0.560"
Set-TextValue "E15" "  -1.68%  "
Set-TextValue "D16" "64.02"
Set-TextValue "E16" "  -2.42%  "
Set-TextValue "D17" "27.532.04"
Set-TextValue "E17" "  +0.61%  "
Set-TextValue "D18" "228.25"
Set-TextValue "E18" "  -1.51%  "
Set-TextValue "E19" "  -0.32%  "
Set-TextValue "D20" "7.61"
Set-TextValue "E20" "  +1.56%  "
Set-TextValue "E21" "  +0.07%  "
Set-TextValue "D22" "4.30"
Set-TextValue "E22" "  -1.36%  "
Set-TextValue "E24" "  -3.26%  "
Set-TextValue "D25" "149.16"
Set-TextValue "E25" "  +1.44%  "
Set-TextValue "D26" "6.95"
Set-TextValue "E26" "  -2.84%  "
Set-TextValue "D27" "0.113"
Set-TextValue "E27" "  +1.41%  "
Set-TextValue "E28" "  -0.08%  "
Set-TextValue "D29" "15.58"
Set-TextValue "E29" "  -1.56%  "
Set-TextValue "E30" "  -0.84%  "
Set-TextValue "D31" "0.0485"
Set-TextValue "E31" "  -2.49%  "
Set-TextValue "E32" "  -0.24%  "
Set-TextValue "D34" "1.426.78"
Set-TextValue "E34" "  -2.28%  "
Set-TextValue "E35" "  +2.23%  "
Set-TextValue "E36" "  -2.03%  "
Set-TextValue "E37" "  +0.24%  "
Set-TextValue "D38" "0.876"
Set-TextValue "E38" "  -3.26%  "
Set-TextValue "E39" "  -1.37%  "
Set-TextValue "D40" "0.910"
Set-TextValue "E40" "  +16.06%  "
Set-TextValue "E41" "  -2.22%  "
Set-TextValue "E42" "  +0.01%  "
Set-TextValue "D44" "5.50"
Set-TextValue "E44" "  +0.99%  "
Set-TextValue "E45" "  +1.89%  "
Set-TextValue "D46" "65.05"
Set-TextValue "E46" "  -0.01%  "
Set-TextValue "D47" "1.782.03"
Set-TextValue "E47" "  -0.69%  "
Set-TextValue "E48" "  -2.40%  "
Set-TextValue "D49" "86.24"
Set-TextValue "E49" "  -2.13%  "
Set-TextValue "E50" "  +0.66%  "
Set-TextValue "D51" "0.0983"
Set-TextValue "E51" "  -2.70%  "
